$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-11-13 Thursday", $false, $false, $false, $false, $false, $true, 1, $false, "2025-11-14 Friday", 2) | Out-Null
$d.Content.Find.Execute("56÷5=", $false, $false, $false, $false, $false, $true, 1, $false, "24÷7=", 2) | Out-Null
$d.Content.Find.Execute("68÷6=", $false, $false, $false, $false, $false, $true, 1, $false, "27÷8=", 2) | Out-Null
$d.Content.Find.Execute("46÷9=", $false, $false, $false, $false, $false, $true, 1, $false, "53÷2=", 2) | Out-Null
$d.Content.Find.Execute("83÷4=", $false, $false, $false, $false, $false, $true, 1, $false, "60÷6=", 2) | Out-Null
$d.Content.Find.Execute("95÷6=", $false, $false, $false, $false, $false, $true, 1, $false, "81÷9=", 2) | Out-Null
$d.Content.Find.Execute("87÷9=", $false, $false, $false, $false, $false, $true, 1, $false, "92÷3=", 2) | Out-Null
$d.Content.Find.Execute("91÷4=", $false, $false, $false, $false, $false, $true, 1, $false, "93÷9=", 2) | Out-Null
$d.Content.Find.Execute("49÷6=", $false, $false, $false, $false, $false, $true, 1, $false, "60÷8=", 2) | Out-Null
$d.Content.Find.Execute("61÷7=", $false, $false, $false, $false, $false, $true, 1, $false, "15÷9=", 2) | Out-Null
$d.Content.Find.Execute("67÷6=", $false, $false, $false, $false, $false, $true, 1, $false, "71÷5=", 2) | Out-Null
$d.Content.Find.Execute("22÷9=", $false, $false, $false, $false, $false, $true, 1, $false, "54÷2=", 2) | Out-Null
$d.Content.Find.Execute("68÷7=", $false, $false, $false, $false, $false, $true, 1, $false, "34÷6=", 2) | Out-Null
$d.Content.Find.Execute("24÷2=", $false, $false, $false, $false, $false, $true, 1, $false, "22÷9=", 2) | Out-Null
$d.Content.Find.Execute("65÷2=", $false, $false, $false, $false, $false, $true, 1, $false, "70÷3=", 2) | Out-Null
$d.Content.Find.Execute("32÷8=", $false, $false, $false, $false, $false, $true, 1, $false, "38÷4=", 2) | Out-Null
$d.Content.Find.Execute("89÷7=", $false, $false, $false, $false, $false, $true, 1, $false, "77÷7=", 2) | Out-Null
$d.Content.Find.Execute("72÷9=", $false, $false, $false, $false, $false, $true, 1, $false, "32÷4=", 2) | Out-Null
$d.Content.Find.Execute("34÷2=", $false, $false, $false, $false, $false, $true, 1, $false, "10÷9=", 2) | Out-Null
$d.Content.Find.Execute("79÷2=", $false, $false, $false, $false, $false, $true, 1, $false, "60÷6=", 2) | Out-Null
$d.Content.Find.Execute("81÷5=", $false, $false, $false, $false, $false, $true, 1, $false, "74÷4=", 2) | Out-Null
$d.Content.Find.Execute("26÷5=", $false, $false, $false, $false, $false, $true, 1, $false, "36÷9=", 2) | Out-Null
$d.Content.Find.Execute("80÷4=", $false, $false, $false, $false, $false, $true, 1, $false, "73÷9=", 2) | Out-Null
$d.Content.Find.Execute("98÷2=", $false, $false, $false, $false, $false, $true, 1, $false, "31÷2=", 2) | Out-Null
$d.Content.Find.Execute("28÷7=", $false, $false, $false, $false, $false, $true, 1, $false, "11÷2=", 2) | Out-Null
$d.Content.Find.Execute("25÷6=", $false, $false, $false, $false, $false, $true, 1, $false, "72÷5=", 2) | Out-Null
